$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 293
$ws.Range("C293").NumberFormat = "@"
$ws.Range("C293").Value = "3"
$ws.Range("H293").Value = 0.028512636263294144
$ws.Range("I293").Value = 0.0006302160761417243
$ws.Range("J293").Value = 2.8512636263294144
$ws.Range("K293").Value = 0.06302160761417243
$ws.Range("L293").Value = 45.242635570093206

# Row 294
$ws.Range("C294").NumberFormat = "@"
$ws.Range("C294").Value = "3"
$ws.Range("H294").Value = 0.04034998083055874
$ws.Range("I294").Value = 0.003695737054460244
$ws.Range("J294").Value = 4.034998083055874
$ws.Range("K294").Value = 0.3695737054460244
$ws.Range("L294").Value = 10.917979346464026

# Row 295
$ws.Range("C295").NumberFormat = "@"
$ws.Range("C295").Value = "3"
$ws.Range("H295").Value = 0.023127654683995136
$ws.Range("I295").Value = 0.0006368675099184532
$ws.Range("J295").Value = 2.3127654683995136
$ws.Range("K295").Value = 0.06368675099184533
$ws.Range("L295").Value = 36.3147033312415

# Row 296
$ws.Range("C296").NumberFormat = "@"
$ws.Range("C296").Value = "3"
$ws.Range("H296").Value = 0.0029910447317695876
$ws.Range("I296").Value = 0.005301571429041399
$ws.Range("J296").Value = 0.29910447317695876
$ws.Range("K296").Value = 0.53015714290414
$ws.Range("L296").Value = 0.5641807852262422

# Row 297
$ws.Range("C297").NumberFormat = "@"
$ws.Range("C297").Value = "3"
$ws.Range("H297").Value = 0.030085769085459013
$ws.Range("I297").Value = 0.0009424389470132288
$ws.Range("J297").Value = 3.0085769085459013
$ws.Range("K297").Value = 0.09424389470132288
$ws.Range("L297").Value = 31.923308327618074

# Row 298
$ws.Range("C298").NumberFormat = "@"
$ws.Range("C298").Value = "3"
$ws.Range("H298").Value = -0.0036801938790174127
$ws.Range("I298").Value = 0.005037006097679386
$ws.Range("J298").Value = -0.36801938790174127
$ws.Range("K298").Value = 0.5037006097679386
$ws.Range("L298").Value = -0.7306312137904549

# Row 299
$ws.Range("C299").NumberFormat = "@"
$ws.Range("C299").Value = "3"
$ws.Range("H299").Value = 0.02756644195157376
$ws.Range("I299").Value = 0.0009470658451166578
$ws.Range("J299").Value = 2.756644195157376
$ws.Range("K299").Value = 0.09470658451166579
$ws.Range("L299").Value = 29.107207374982654

# Row 300
$ws.Range("C300").NumberFormat = "@"
$ws.Range("C300").Value = "3"
$ws.Range("H300").Value = -0.031302985901833114
$ws.Range("I300").Value = 0.00781493740746596
$ws.Range("J300").Value = -3.1302985901833114
$ws.Range("K300").Value = 0.7814937407465959
$ws.Range("L300").Value = -4.005532516732376

# Row 301
$ws.Range("C301").NumberFormat = "@"
$ws.Range("C301").Value = "3"
$ws.Range("H301").Value = 0.01704127853410342
$ws.Range("I301").Value = 0.000966769276173102
$ws.Range("J301").Value = 1.7041278534103421
$ws.Range("K301").Value = 0.0966769276173102
$ws.Range("L301").Value = 17.627037757716398

# Row 302
$ws.Range("C302").NumberFormat = "@"
$ws.Range("C302").Value = "3"
$ws.Range("H302").Value = 0.06911330739768662
$ws.Range("I302").Value = 0.007290734502166438
$ws.Range("J302").Value = 6.911330739768662
$ws.Range("K302").Value = 0.7290734502166438
$ws.Range("L302").Value = 9.479608313421595

# Row 303
$ws.Range("C303").NumberFormat = "@"
$ws.Range("C303").Value = "3"
$ws.Range("H303").Value = 0.012184273858861783
$ws.Range("I303").Value = 0.000976069695138729
$ws.Range("J303").Value = 1.2184273858861783
$ws.Range("K303").Value = 0.0976069695138729
$ws.Range("L303").Value = 12.48299575280844

# Row 304
$ws.Range("C304").NumberFormat = "@"
$ws.Range("C304").Value = "3"
$ws.Range("H304").Value = 0.07144126969077313
$ws.Range("I304").Value = 0.010162722070237684
$ws.Range("J304").Value = 7.1441269690773135
$ws.Range("K304").Value = 1.0162722070237684
$ws.Range("L304").Value = 7.029737623150632

# Row 305
$ws.Range("B305").Value = 1
$ws.Range("C305").NumberFormat = "@"
$ws.Range("C305").Value = "6"
$ws.Range("E305").Value = 0.004
$ws.Range("G305").Value = 999
$ws.Range("H305").Value = 0.015518657789766799
$ws.Range("I305").Value = 0.0006172615230912757
$ws.Range("J305").Value = 1.5518657789766799
$ws.Range("K305").Value = 0.06172615230912757
$ws.Range("L305").Value = 25.141139062173526

# Row 306
$ws.Range("B306").Value = 2
$ws.Range("C306").NumberFormat = "@"
$ws.Range("C306").Value = "2"
$ws.Range("E306").Value = 0.004
$ws.Range("G306").Value = 0.004
$ws.Range("H306").Value = 0.019901956072249094
$ws.Range("I306").Value = 0.0027732343367971994
$ws.Range("J306").Value = 1.9901956072249094
$ws.Range("K306").Value = 0.27732343367971996
$ws.Range("L306").Value = 7.176442253067517

# Row 307
$ws.Range("B307").Value = 2
$ws.Range("C307").NumberFormat = "@"
$ws.Range("C307").Value = "4"
$ws.Range("E307").Value = 0.005
$ws.Range("F307").Value = 0.0968
$ws.Range("G307").Value = 0.004
$ws.Range("H307").Value = 0.04051999278931673
$ws.Range("I307").Value = 0.0014209560623960771
$ws.Range("J307").Value = 4.051999278931673
$ws.Range("K307").Value = 0.1420956062396077
$ws.Range("L307").Value = 28.51600683626359

# Row 308
$ws.Range("B308").Value = 1
$ws.Range("C308").NumberFormat = "@"
$ws.Range("C308").Value = "6"
$ws.Range("E308").Value = 0.004
$ws.Range("G308").Value = 0.005
$ws.Range("H308").Value = 0.016473167423109114
$ws.Range("I308").Value = 0.0006143687926401384
$ws.Range("J308").Value = 1.6473167423109114
$ws.Range("K308").Value = 0.061436879264013845
$ws.Range("L308").Value = 26.81315786291596

# Row 309
$ws.Range("B309").Value = 2
$ws.Range("C309").NumberFormat = "@"
$ws.Range("C309").Value = "2"
$ws.Range("E309").Value = 0.004
$ws.Range("G309").Value = 0.004
$ws.Range("H309").Value = 0.02078401241398753
$ws.Range("I309").Value = 0.0027708379934922975
$ws.Range("J309").Value = 2.078401241398753
$ws.Range("K309").Value = 0.2770837993492298
$ws.Range("L309").Value = 7.500984345819461

# Row 310
$ws.Range("B310").Value = 2
$ws.Range("C310").NumberFormat = "@"
$ws.Range("C310").Value = "4"
$ws.Range("E310").Value = 0.005
$ws.Range("F310").Value = 0.103
$ws.Range("G310").Value = 0.004
$ws.Range("H310").Value = 0.04246738028188335
$ws.Range("I310").Value = 0.0014130076522895964
$ws.Range("J310").Value = 4.246738028188335
$ws.Range("K310").Value = 0.14130076522895965
$ws.Range("L310").Value = 30.05460035058582

# Row 311
$ws.Range("B311").Value = 1
$ws.Range("C311").NumberFormat = "@"
$ws.Range("C311").Value = "6"
$ws.Range("E311").Value = 0.004
$ws.Range("G311").Value = 0.005
$ws.Range("H311").Value = 0.01591948205915439
$ws.Range("I311").Value = 0.0006160448014426987
$ws.Range("J311").Value = 1.5919482059154388
$ws.Range("K311").Value = 0.06160448014426987
$ws.Range("L311").Value = 25.841435593438955

# Row 312
$ws.Range("B312").Value = 2
$ws.Range("C312").NumberFormat = "@"
$ws.Range("C312").Value = "2"
$ws.Range("E312").Value = 0.004
$ws.Range("G312").Value = 0.004
$ws.Range("H312").Value = 0.01960776772247086
$ws.Range("I312").Value = 0.002774034500604222
$ws.Range("J312").Value = 1.960776772247086
$ws.Range("K312").Value = 0.2774034500604222
$ws.Range("L312").Value = 7.068321507248747

# Row 313
$ws.Range("B313").Value = 2
$ws.Range("C313").NumberFormat = "@"
$ws.Range("C313").Value = "4"
$ws.Range("E313").Value = 0.005
$ws.Range("F313").Value = 0.09939999999999999
$ws.Range("G313").Value = 0.004
$ws.Range("H313").Value = 0.04171627485608709
$ws.Range("I313").Value = 0.0014160663060587787
$ws.Range("J313").Value = 4.171627485608709
$ws.Range("K313").Value = 0.14160663060587786
$ws.Range("L313").Value = 29.459266615976887

# Row 314
$ws.Range("B314").Value = 1
$ws.Range("C314").NumberFormat = "@"
$ws.Range("C314").Value = "6"
$ws.Range("E314").Value = 0.01
$ws.Range("G314").Value = 0.005
$ws.Range("H314").Value = 0.008963393925290664
$ws.Range("I314").Value = 0.0015939390109404278
$ws.Range("J314").Value = 0.8963393925290664
$ws.Range("K314").Value = 0.15939390109404278
$ws.Range("L314").Value = 5.6234233956054815

# Row 315
$ws.Range("B315").Value = 2
$ws.Range("C315").NumberFormat = "@"
$ws.Range("C315").Value = "2"
$ws.Range("E315").Value = 0.012
$ws.Range("G315").Value = 0.01
$ws.Range("H315").Value = 0.027131929208706884
$ws.Range("I315").Value = 0.00760394011110491
$ws.Range("J315").Value = 2.7131929208706884
$ws.Range("K315").Value = 0.760394011110491
$ws.Range("L315").Value = 3.568140833866248

# Row 316
$ws.Range("B316").Value = 2
$ws.Range("C316").NumberFormat = "@"
$ws.Range("C316").Value = "4"
$ws.Range("E316").Value = 0.017
$ws.Range("F316").Value = 0.055
$ws.Range("G316").Value = 0.01
$ws.Range("H316").Value = 0.03869548797561162
$ws.Range("I316").Value = 0.004399973634324998
$ws.Range("J316").Value = 3.869548797561162
$ws.Range("K316").Value = 0.4399973634324998
$ws.Range("L316").Value = 8.794481783650033

# Row 317
$ws.Range("B317").Value = 1
$ws.Range("C317").NumberFormat = "@"
$ws.Range("C317").Value = "6"
$ws.Range("E317").Value = 0.01
$ws.Range("G317").Value = 0.017
$ws.Range("H317").Value = 0.009917498912697331
$ws.Range("I317").Value = 0.0015864239693884658
$ws.Range("J317").Value = 0.9917498912697331
$ws.Range("K317").Value = 0.15864239693884657
$ws.Range("L317").Value = 6.251480754240196

# Row 318
$ws.Range("B318").Value = 2
$ws.Range("C318").NumberFormat = "@"
$ws.Range("C318").Value = "2"
$ws.Range("E318").Value = 0.012
$ws.Range("G318").Value = 0.01
$ws.Range("H318").Value = 0.030048542545447088
$ws.Range("I318").Value = 0.007582409326657588
$ws.Range("J318").Value = 3.004854254544709
$ws.Range("K318").Value = 0.7582409326657589
$ws.Range("L318").Value = 3.962928041856694

# Row 319
$ws.Range("B319").Value = 2
$ws.Range("C319").NumberFormat = "@"
$ws.Range("C319").Value = "4"
$ws.Range("E319").Value = 0.017
$ws.Range("F319").Value = 0.061
$ws.Range("G319").Value = 0.01
$ws.Range("H319").Value = 0.04334896024957158
$ws.Range("I319").Value = 0.004341362460117689
$ws.Range("J319").Value = 4.334896024957158
$ws.Range("K319").Value = 0.4341362460117689
$ws.Range("L319").Value = 9.985105055797723

# Row 320
$ws.Range("B320").Value = 1
$ws.Range("C320").NumberFormat = "@"
$ws.Range("C320").Value = "6"
$ws.Range("E320").Value = 0.009000000000000001
$ws.Range("G320").Value = 0.017
$ws.Range("H320").Value = 0.00784453941396146
$ws.Range("I320").Value = 0.0014425255812222733
$ws.Range("J320").Value = 0.784453941396146
$ws.Range("K320").Value = 0.14425255812222731
$ws.Range("L320").Value = 5.438059134670365

# Row 321
$ws.Range("B321").Value = 2
$ws.Range("C321").NumberFormat = "@"
$ws.Range("C321").Value = "2"
$ws.Range("E321").Value = 0.011
$ws.Range("G321").Value = 0.009000000000000001
$ws.Range("H321").Value = 0.02956301409869999
$ws.Range("I321").Value = 0.006902282914656734
$ws.Range("J321").Value = 2.956301409869999
$ws.Range("K321").Value = 0.6902282914656734
$ws.Range("L321").Value = 4.283077709829028

# Row 322
$ws.Range("B322").Value = 2
$ws.Range("C322").NumberFormat = "@"
$ws.Range("C322").Value = "4"
$ws.Range("E322").Value = 0.016
$ws.Range("F322").Value = 0.048
$ws.Range("G322").Value = 0.009000000000000001
$ws.Range("H322").Value = 0.04136225058241694
$ws.Range("I322").Value = 0.004063960403311092
$ws.Range("J322").Value = 4.136225058241694
$ws.Range("K322").Value = 0.4063960403311092
$ws.Range("L322").Value = 10.177818304705243

# Row 323
$ws.Range("B323").Value = 1
$ws.Range("C323").NumberFormat = "@"
$ws.Range("C323").Value = "6"
$ws.Range("E323").Value = 0.011
$ws.Range("G323").Value = 0.016
$ws.Range("H323").Value = 0.008644354023879863
$ws.Range("I323").Value = 0.0017561076122606326
$ws.Range("J323").Value = 0.8644354023879863
$ws.Range("K323").Value = 0.17561076122606326
$ws.Range("L323").Value = 4.922451200329352

# Row 324
$ws.Range("B324").Value = 2
$ws.Range("C324").NumberFormat = "@"
$ws.Range("C324").Value = "2"
$ws.Range("E324").Value = 0.013
$ws.Range("G324").Value = 0.011
$ws.Range("H324").Value = 0.03295692068933831
$ws.Range("I324").Value = 0.008243028351347862
$ws.Range("J324").Value = 3.295692068933831
$ws.Range("K324").Value = 0.8243028351347862
$ws.Range("L324").Value = 3.9981569011526386

# Row 325
$ws.Range("B325").Value = 2
$ws.Range("C325").NumberFormat = "@"
$ws.Range("C325").Value = "4"
$ws.Range("E325").Value = 0.018
$ws.Range("F325").Value = 0.053
$ws.Range("G325").Value = 0.011
$ws.Range("H325").Value = 0.04069748021458408
$ws.Range("I325").Value = 0.004678929538653746
$ws.Range("J325").Value = 4.069748021458408
$ws.Range("K325").Value = 0.4678929538653746
$ws.Range("L325").Value = 8.698032291012838
